# Localization status report regeneration:
#  - Source strings moved from "Ready for handoff" to "In Translation"
#    for both locales (zh-cn, de-de).
#  - Narrower "Status" columns follow the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhCn = $wb.Worksheets.Item("zh-cn")
$deDe = $wb.Worksheets.Item("de-de")

# Update the localization status value everywhere it is shown: the two
# locale detail sheets, and the corresponding columns on the Overview
# summary sheet.
$zhCn.Range("C2").Value = "In Translation"
$deDe.Range("C2").Value = "In Translation"
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Shrink the now-narrower "Status" columns to match the regenerated report.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhCn.Columns.Item(3).ColumnWidth = 12.5
$deDe.Columns.Item(3).ColumnWidth = 12.5
